$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a run split at a single (collapsed) point by briefly adding
# and then removing a bookmark there. The underlying engine never re-merges
# runs that are already split in the in-memory model, so the split survives
# the bookmark's removal -- this lets us reproduce Word's "typing splits a
# run" behaviour from a script.
# ---------------------------------------------------------------------------
function Split-At($pos) {
    $tmpName = "zzTmpSplit"
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($tmpName, $r)
    $d.Bookmarks($tmpName).Delete()
}

# ===========================================================================
# 1) Title: "neonCLUSTER Node Templates" -> "neonHIVE Node Templates"
#    Split into separate "neon" / "HIVE" / " Node " / "Templates" runs.
#    Drop the _GoBack bookmark that used to wrap "neonCLUSTER" here -- it
#    relocates to the download-link URL below.
# ===========================================================================
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$titleFind = $d.Content
$titleFind.Find.Execute("neonCLUSTER", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$titleStart = $titleFind.Start

$clusterRange = $d.Range($titleStart + 4, $titleStart + 11)   # "CLUSTER"
$clusterRange.Text = "HIVE"

# Do every text substitution for this region *before* splitting -- setting
# .Text re-merges neighbouring same-format runs, so any split made earlier
# would otherwise be undone.
Split-At ($titleStart + 14)    # before "Templates"
Split-At ($titleStart + 8)     # before " Node "
Split-At ($titleStart + 4)     # before "HIVE"

# ===========================================================================
# 2) Hyperlink URL (first occurrence only):
#    https://s3.amazonaws.com/neonforge/neoncluster/neon-ubuntu-16.04.#.vhdx
#    Text is unchanged, but it now gets split into 3 runs, with the
#    relocated _GoBack bookmark wrapping just "cluster".
# ===========================================================================
$urlFind = $d.Content
$urlFind.Find.Execute("https://s3.amazonaws.com/neonforge/neoncluster/neon-ubuntu-16.04.#.vhdx", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$urlStart = $urlFind.Start

$clusterStart = $urlStart + ("https://s3.amazonaws.com/neonforge/neon").Length
$clusterEnd = $clusterStart + ("cluster").Length

Split-At $clusterEnd
Split-At $clusterStart
$bmRange = $d.Range($clusterStart, $clusterEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ===========================================================================
# 3) "Create a 10GB VHD (the disk will be resized during cluster
#     provisioning)." -> "... during hive provisioning)."
# ===========================================================================
$p17Find = $d.Content
$p17Find.Find.Execute("the disk will be resized during cluster provisioning", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p17Start = $p17Find.Start

$hiveStart = $p17Start + ("the disk will be resized during ").Length
$hiveEnd = $hiveStart + ("cluster").Length
$hiveRange = $d.Range($hiveStart, $hiveEnd)
$hiveRange.Text = "hive"
$hiveEnd = $hiveStart + ("hive").Length

Split-At $hiveEnd
Split-At $hiveStart

# ===========================================================================
# 4) "... A secure password will be set during cluster setup."
#     -> "... will be set during hive setup."
# ===========================================================================
$p30Find = $d.Content
$p30Find.Find.Execute("template credentials.  A secure password will be set during cluster setup.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p30Start = $p30Find.Start

$hive2Start = $p30Start + ("template credentials.  A secure password will be set during ").Length
$hive2End = $hive2Start + ("cluster").Length
$hive2Range = $d.Range($hive2Start, $hive2End)
$hive2Range.Text = "hive"
$hive2End = $hive2Start + ("hive").Length

Split-At $hive2End
Split-At $hive2Start
